$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 9, shifting existing rows 9-18 down to 10-19.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new weekly record.
$ws.Range("A9").Value = 2
$ws.Range("B9").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C9").Value = "Coquimbo"
$ws.Range("D9").Value = 45203
$ws.Range("D9").NumberFormat = $ws.Range("D10").NumberFormat
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100101
$ws.Range("H9").Value = "Berries"
$ws.Range("I9").Value = 100101001
$ws.Range("J9").Value = "Arándano (blue)"
$ws.Range("K9").Value = "Sin especificar"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 1000
$ws.Range("N9").Value = 9000
$ws.Range("O9").Value = 10000
$ws.Range("P9").Value = 9500
$ws.Range("Q9").Value = "$/bandeja 2 kilos"
$ws.Range("R9").Value = "Provincia de Limarí"
$ws.Range("S9").Value = 4750
$ws.Range("T9").Value = 2
